# Update the cryptos list "Price" (D) and "Volume(1h)" (E) columns with the
# latest scraped figures. Price strings that look like plain numbers are
# prefixed with a leading apostrophe so Excel stores them as literal text
# (matching the source data, which keeps e.g. "23.60" / "1.0000" with their
# trailing zeros instead of collapsing them into a numeric value).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.459.51'
$ws.Range("E2").Value = '  +1.77%  '
$ws.Range("D3").Value = '1.858.69'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("D4").Value = "'" + '0.9998'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'" + '244.92'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").Value = "'" + '0.6947'
$ws.Range("E6").Value = '  +0.79%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = "'" + '0.07685'
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = "'" + '0.3063'
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").Value = "'" + '23.60'
$ws.Range("E10").Value = '  +0.32%  '
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").Value = '1.850.57'
$ws.Range("E13").Value = '  +0.84%  '
$ws.Range("D14").Value = "'" + '91.19'
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").Value = "'" + '0.6922'
$ws.Range("D16").Value = "'" + '6.331'
$ws.Range("E16").Value = '  -1.78%  '
$ws.Range("D17").Value = '29.457.52'
$ws.Range("E17").Value = '  +1.60%  '
$ws.Range("D18").Value = "'" + '0.000008304'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '2.101.90'
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").Value = "'" + '238.24'
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = "'" + '1.0000'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = "'" + '7.630'
$ws.Range("E23").Value = '  +2.03%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = "'" + '0.1493'
$ws.Range("E25").Value = '  +1.62%  '
$ws.Range("D26").Value = "'" + '8.901'
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("D27").Value = "'" + '159.81'
$ws.Range("E27").Value = '  -2.13%  '
$ws.Range("D28").Value = "'" + '18.26'
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").Value = "'" + '1.533'
$ws.Range("E29").Value = '  -1.29%  '
$ws.Range("E30").Value = '  +0.93%  '
$ws.Range("D31").Value = "'" + '4.156'
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  +3.65%  '
$ws.Range("D33").Value = "'" + '0.05114'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = "'" + '0.7720'
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").Value = "'" + '1.883'
$ws.Range("E35").Value = '  +2.06%  '
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").Value = "'" + '2.684'
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").Value = '1.329.52'
$ws.Range("E38").Value = '  +7.06%  '
$ws.Range("D39").Value = "'" + '0.01871'
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("D40").Value = "'" + '2.721'
$ws.Range("E40").Value = '  +0.90%  '
$ws.Range("D41").Value = "'" + '0.9549'
$ws.Range("E41").Value = '  +1.31%  '
$ws.Range("D42").Value = "'" + '5.859'
$ws.Range("E42").Value = '  +2.81%  '
$ws.Range("D43").Value = "'" + '105.76'
$ws.Range("E43").Value = '  -2.43%  '
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").Value = "'" + '9.833'
$ws.Range("E45").Value = '  +2.59%  '
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("D47").Value = '2.000.78'
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("D48").Value = "'" + '0.5220'
$ws.Range("E48").Value = '  +0.89%  '
$ws.Range("D49").Value = "'" + '1.784'
$ws.Range("E49").Value = '  +2.03%  '
$ws.Range("D50").Value = "'" + '63.36'
$ws.Range("E50").Value = '  -1.76%  '
$ws.Range("D51").Value = "'" + '6.975'
$ws.Range("E51").Value = '  +0.84%  '
